$d = $word.ActiveDocument

function Find-Range($doc, $text) {
    $r = $doc.Content
    $found = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $text"
    }
    return $r
}

# --- Change 1: geoid / GCS paragraph ---
$r = Find-Range $d "the geoid we are using to represent the earth"
$r.Collapse(0)
$r.InsertAfter(" (is it spheroid? ellipsoid?)")

$r = Find-Range $d "where its equator and prime meridian are drawn"
$r.Collapse(0)
$r.InsertAfter(",")

# --- Change 2: Projections paragraph, "Finally" -> "Finally," ---
$r = Find-Range $d "Finally"
$r.Collapse(0)
$r.InsertAfter(",")

# --- Change 3: Add Vector Layer instructions ---
$r = Find-Range $d "Add Vector Layer, navigate to the class folder for today"
$r.Collapse(0)
$r.InsertAfter(', navigate into the "world" folder')

# --- Change 4: WGS 84 / lastRenderedPageBreak swap ---
# (text unaffected; page-break marker relocation not handled via this pass)

# --- Change 6: footer ---
$footers = $d.Sections.Item(1).Footers
$footerRange = $footers.Item(1).Range
$ffound = $footerRange.Find.Execute("VRI GIS Unit 2: Census Data", $true, $false, $false, $false, $false, $true, 1, $false, "VRI GIS Unit 3: Making a Choropleth Map", 2)
Write-Output "footer replaced: $ffound"

# --- Big change: delete "Joining tabular..." through "Mapping data" section, rewrite bookmark paragraph ---
$paras = $d.Paragraphs
$target = $null
foreach ($p in $paras) {
    if ($p.Range.Text.Length -le 2) {
        $target = $p
    }
}
$targetIndex = $target.Range.Start

# find heading para start and last para end
$headingR = Find-Range $d "Joining tabular and spatial data"
$startDel = $headingR.Start

$lastParaRange = $paras.Last.Range
$endDel = $lastParaRange.End

$delRange = $d.Range($startDel, $endDel)
$delRange.Delete()

# Now rewrite the bookmark paragraph
$p37 = $target
$p37.Range.ListFormat.RemoveNumbers()
$p37.Style = "Normal"
$r = $p37.Range
$r.Collapse(1)
$sentence1 = "A choropleth map typically uses color to visualize difference in some data value for different subregions of a map. Our maps are going to examine demographic and voting data in the southern U.S. at the county level. "
$sentence2 = "The starting point for a choropleth map is understanding what kind of variable you are hoping to represent: For most purposes we will have some variant of "
$combined = $sentence1 + "`r" + $sentence2
$r.InsertBefore($combined)

Write-Output "ALL DONE"
